$d = $word.ActiveDocument

# Locate the end of the paragraph that ends with "...make everything neater."
# so we can insert the new paragraph right after it.
$r = $d.Content
$found = $r.Find.Execute("make everything neater.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos = $r.End

# Build a fresh, collapsed range at that position (re-using the Find range
# directly for InsertXML can clobber neighbouring content), then insert the
# new paragraph's OOXML directly so formatting/proofErr marks match exactly.
$ins = $d.Range($pos, $pos)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="283"/><w:jc w:val="both"/><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-US"/></w:rPr><w:t>So</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> if you want to guarantee at most one solution just by counting clues, you need at least 78.</w:t></w:r></w:p>'

$ins.InsertXML($xml)
